$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.215492537586345
$ws.Range("D2").Value = 4.064204902707336
$ws.Range("E2").Value = 16.53398874325575
$ws.Range("F2").Value = 22.02924588375378
$ws.Range("G2").Value = 3.606149635419113
$ws.Range("I2").Value = 18.33330853334726
$ws.Range("K2").Value = 10.95075768655988
$ws.Range("N2").Value = 17.16838657490608
$ws.Range("O2").Value = 19.50962468359262
$ws.Range("B3").Value = 6.097650727729514
$ws.Range("D3").Value = 4.027265270322727
$ws.Range("E3").Value = 15.59360508834405
$ws.Range("F3").Value = 21.95246083352328
$ws.Range("G3").Value = 3.60830338140574
$ws.Range("I3").Value = 18.43199165608033
$ws.Range("K3").Value = 10.40684809081244
$ws.Range("N3").Value = 17.22190143745738
$ws.Range("O3").Value = 19.51857935414288
$ws.Range("B4").Value = 6.025040847459481
$ws.Range("D4").Value = 4.004080442784298
$ws.Range("E4").Value = 14.99119597541937
$ws.Range("F4").Value = 21.91259646262687
$ws.Range("G4").Value = 3.609694608971342
$ws.Range("I4").Value = 18.49563387082981
$ws.Range("K4").Value = 10.05583208996558
$ws.Range("N4").Value = 17.25644739765533
$ws.Range("O4").Value = 19.52980250078911
$ws.Range("B5").Value = 5.995427612837357
$ws.Range("D5").Value = 3.994510051356479
$ws.Range("E5").Value = 14.73968927372237
$ws.Range("F5").Value = 21.89819316180543
$ws.Range("G5").Value = 3.610278910209348
$ws.Range("I5").Value = 18.52233774066834
$ws.Range("K5").Value = 9.908591338360681
$ws.Range("N5").Value = 17.27095068193547
$ws.Range("O5").Value = 19.53581179088993
$ws.Range("B6").Value = 5.990510083995349
$ws.Range("D6").Value = 3.992913629774495
$ws.Range("E6").Value = 14.69757208222217
$ws.Range("F6").Value = 21.89591300018574
$ws.Range("G6").Value = 3.610376983524037
$ws.Range("I6").Value = 18.52681841477794
$ws.Range("K6").Value = 9.883891469188551
$ws.Range("N6").Value = 17.27338467615893
$ws.Range("O6").Value = 19.53689621624626
$ws.Range("B7").Value = 6.024641517314153
$ws.Range("D7").Value = 4.003951862989574
$ws.Range("E7").Value = 14.98782804902727
$ws.Range("F7").Value = 21.91239474550287
$ws.Range("G7").Value = 3.609702418674287
$ws.Range("I7").Value = 18.49599089150456
$ws.Range("K7").Value = 10.05386322029593
$ws.Range("N7").Value = 17.25664126955509
$ws.Range("O7").Value = 19.52987773660884
$ws.Range("B8").Value = 6.174937143137135
$ws.Range("D8").Value = 4.051575216457238
$ws.Range("E8").Value = 16.21508082142385
$ws.Range("F8").Value = 22.00126622920677
$ws.Range("G8").Value = 3.606877999673726
$ws.Range("I8").Value = 18.36670264138089
$ws.Range("K8").Value = 10.76680556501608
$ws.Range("N8").Value = 17.18648893259389
$ws.Range("O8").Value = 19.51152201757706
$ws.Range("B9").Value = 6.465947717850471
$ws.Range("D9").Value = 4.140784035996422
$ws.Range("E9").Value = 18.51239503973919
$ws.Range("F9").Value = 22.23272834179156
$ws.Range("G9").Value = 3.60188264383834
$ws.Range("I9").Value = 18.13727489994499
$ws.Range("K9").Value = 12.02644235361654
$ws.Range("N9").Value = 17.06225786094262
$ws.Range("O9").Value = 19.52108195871286
$ws.Range("B10").Value = 6.675304082294837
$ws.Range("D10").Value = 4.20353876552366
$ws.Range("E10").Value = 20.15899490220449
$ws.Range("F10").Value = 22.43667276882572
$ws.Range("G10").Value = 3.598539959742358
$ws.Range("I10").Value = 17.98327244270785
$ws.Range("K10").Value = 12.86423670263255
$ws.Range("N10").Value = 16.97904173080366
$ws.Range("O10").Value = 19.55599683003284
$ws.Range("B11").Value = 6.769149396379895
$ws.Range("D11").Value = 4.231433579599829
$ws.Range("E11").Value = 20.86580858018392
$ws.Range("F11").Value = 22.53655120927201
$ws.Range("G11").Value = 3.597089564258371
$ws.Range("I11").Value = 17.91634460326802
$ws.Range("K11").Value = 13.2258628014715
$ws.Range("N11").Value = 16.94291824867474
$ws.Range("O11").Value = 19.5779415094599
$ws.Range("B12").Value = 6.80445260873687
$ws.Range("D12").Value = 4.241898770508691
$ws.Range("E12").Value = 21.12742418035783
$ws.Range("F12").Value = 22.57536851996653
$ws.Range("G12").Value = 3.596550371013197
$ws.Range("I12").Value = 17.89144849912289
$ws.Range("K12").Value = 13.35996951501159
$ws.Range("N12").Value = 16.92948712518849
$ws.Range("O12").Value = 19.58712164409046
$ws.Range("B13").Value = 6.796860353536592
$ws.Range("D13").Value = 4.239649326053041
$ws.Range("E13").Value = 21.071348650697
$ws.Range("F13").Value = 22.566964677078
$ws.Range("G13").Value = 3.596666050277767
$ws.Range("I13").Value = 17.89679042468549
$ws.Range("K13").Value = 13.33121365681836
$ws.Range("N13").Value = 16.9323687413188
$ws.Range("O13").Value = 19.58510587251786
$ws.Range("B14").Value = 6.772058688674691
$ws.Range("D14").Value = 4.232296543823728
$ws.Range("E14").Value = 20.88745269977489
$ws.Range("F14").Value = 22.53972493722998
$ws.Range("G14").Value = 3.597045003607128
$ws.Range("I14").Value = 17.91428741758501
$ws.Range("K14").Value = 13.23695277413288
$ws.Range("N14").Value = 16.94180829631951
$ws.Range("O14").Value = 19.57867933776465
$ws.Range("B15").Value = 6.756835513093054
$ws.Range("D15").Value = 4.22777987342916
$ws.Range("E15").Value = 20.77402560618524
$ws.Range("F15").Value = 22.52316864964315
$ws.Range("G15").Value = 3.597278429322744
$ws.Range("I15").Value = 17.92506312041461
$ws.Range("K15").Value = 13.17884550305826
$ws.Range("N15").Value = 16.94762257111765
$ws.Range("O15").Value = 19.57485615736659
$ws.Range("B16").Value = 6.669140125628703
$ws.Range("D16").Value = 4.201702286282367
$ws.Range("E16").Value = 20.11195601327193
$ws.Range("F16").Value = 22.4302862457624
$ws.Range("G16").Value = 3.598636154488881
$ws.Range("I16").Value = 17.98770912970222
$ws.Range("K16").Value = 12.84020861484655
$ws.Range("N16").Value = 16.98143725097023
$ws.Range("O16").Value = 19.55468456642986
$ws.Range("B17").Value = 6.614960656856948
$ws.Range("D17").Value = 4.185534408762158
$ws.Range("E17").Value = 19.69500491714664
$ws.Range("F17").Value = 22.37510744102739
$ws.Range("G17").Value = 3.599487017107873
$ws.Range("I17").Value = 18.02694040586605
$ws.Range("K17").Value = 12.62744982575251
$ws.Range("N17").Value = 17.00262434128723
$ws.Range("O17").Value = 19.54386176726349
$ws.Range("B18").Value = 6.583668942755646
$ws.Range("D18").Value = 4.176173884442648
$ws.Range("E18").Value = 19.45120852203443
$ws.Range("F18").Value = 22.34404032358751
$ws.Range("G18").Value = 3.599983022430456
$ws.Range("I18").Value = 18.04979979279821
$ws.Range("K18").Value = 12.50324336625255
$ws.Range("N18").Value = 17.0149736733358
$ws.Range("O18").Value = 19.53820749561705
$ws.Range("B19").Value = 6.573052992484478
$ws.Range("D19").Value = 4.17299418999562
$ws.Range("E19").Value = 19.36797914319819
$ws.Range("F19").Value = 22.33363740469738
$ws.Range("G19").Value = 3.600152098574467
$ws.Range("I19").Value = 18.05759024295731
$ws.Range("K19").Value = 12.46087537975127
$ws.Range("N19").Value = 17.0191829844923
$ws.Range("O19").Value = 19.53639110607866
$ws.Range("B20").Value = 6.620741768063812
$ws.Range("D20").Value = 4.187261869940768
$ws.Range("E20").Value = 19.73980144921946
$ws.Range("F20").Value = 22.38091212906565
$ws.Range("G20").Value = 3.599395757453215
$ws.Range("I20").Value = 18.02273369378362
$ws.Range("K20").Value = 12.65028833766728
$ws.Range("N20").Value = 17.00035206864645
$ws.Range("O20").Value = 19.54495481321228
$ws.Range("B21").Value = 6.779350148409524
$ws.Range("D21").Value = 4.234458925358024
$ws.Range("E21").Value = 20.94163103890711
$ws.Range("F21").Value = 22.5476991082673
$ws.Range("G21").Value = 3.596933423787688
$ws.Range("I21").Value = 17.90913598249226
$ws.Range("K21").Value = 13.26471658260845
$ws.Range("N21").Value = 16.93902894591341
$ws.Range("O21").Value = 19.58054336817629
$ws.Range("B22").Value = 6.881631853337661
$ws.Range("D22").Value = 4.264732056103289
$ws.Range("E22").Value = 21.69190935595179
$ws.Range("F22").Value = 22.66249315943162
$ws.Range("G22").Value = 3.59538264179451
$ws.Range("I22").Value = 17.837503833814
$ws.Range("K22").Value = 13.64975714334403
$ws.Range("N22").Value = 16.90039614520993
$ws.Range("O22").Value = 19.6088730594827
$ws.Range("B23").Value = 6.827178929686176
$ws.Range("D23").Value = 4.248628466653047
$ws.Range("E23").Value = 21.29468023737045
$ws.Range("F23").Value = 22.60070481305698
$ws.Range("G23").Value = 3.596204989356694
$ws.Range("I23").Value = 17.87549701274947
$ws.Range("K23").Value = 13.44577402554357
$ws.Range("N23").Value = 16.92088326849718
$ws.Range("O23").Value = 19.59328979076383
$ws.Range("B24").Value = 6.618128571192631
$ws.Range("D24").Value = 4.186481088093973
$ws.Range("E24").Value = 19.719561666235
$ws.Range("F24").Value = 22.37828578443331
$ws.Range("G24").Value = 3.599436994653037
$ws.Range("I24").Value = 18.02463459870427
$ws.Range("K24").Value = 12.63996892672794
$ws.Range("N24").Value = 17.00137883784714
$ws.Range("O24").Value = 19.54445887795617
$ws.Range("B25").Value = 6.387851014552985
$ws.Range("D25").Value = 4.117123186935047
$ws.Range("E25").Value = 17.86834494980866
$ws.Range("F25").Value = 22.16407985233093
$ws.Range("G25").Value = 3.603176249711677
$ws.Range("I25").Value = 18.19677475506023
$ws.Range("K25").Value = 11.70087971030472
$ws.Range("N25").Value = 17.09444554941775
$ws.Range("O25").Value = 19.52987773660884
